$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update the reported quarter's date range (row 8) ---
# B8: Fecha de inicio del periodo que se informa (period start)
# C8: Fecha de término del periodo que se informa (period end)
# F8: Fecha de validación (validation date)
# G8: Fecha de actualización (update date)
$ws.Range("B8").Value = 44470
$ws.Range("C8").Value = 44561
$ws.Range("F8").Value = 44571
$ws.Range("G8").Value = 44571

# --- Update the saved window view (scroll position + active selection) ---
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1
$ws.Range("C12").Select()
